# Auto-generated edit script: applies the scheduled-runner market-data refresh
# to the per-sheet Leve profit columns (H..N) as captured in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 182.33333
$ws.Range("I2").Value = 65
$ws.Range("J2").Value = 299.66666
$ws.Range("K2").Value = 65
$ws.Range("L2").Value = 299.66666
$ws.Range("M2").Value = 48
$ws.Range("N2").Value = -525.66666
$ws.Range("H18").Value = 1439.2
$ws.Range("I18").Value = 1488.1111
$ws.Range("K18").Value = 1488.1111
$ws.Range("M18").Value = -1204.1111
$ws.Range("H33").Value = 1015.7273
$ws.Range("I33").Value = 1015.7273
$ws.Range("K33").Value = 1015.7273
$ws.Range("M33").Value = -786.7273
$ws.Range("H38").Value = 2496.125
$ws.Range("I38").Value = 61.5
$ws.Range("K38").Value = 184.5
$ws.Range("M38").Value = 187.5
$ws.Range("H42").Value = 28
$ws.Range("I42").Value = 28
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 84
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 146
$ws.Range("N42").ClearContents()
$ws.Range("H51").Value = 27999.6
$ws.Range("I51").Value = 62499
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 62499
$ws.Range("L51").Value = 5000
$ws.Range("M51").Value = -62015
$ws.Range("N51").Value = -5968
$ws.Range("H53").Value = 1181.4
$ws.Range("I53").Value = 1226.75
$ws.Range("K53").Value = 1226.75
$ws.Range("M53").Value = -589.75
$ws.Range("H64").Value = 4166.5
$ws.Range("I64").Value = 3999.8
$ws.Range("K64").Value = 3999.8
$ws.Range("M64").Value = -3751.8
$ws.Range("H67").Value = 4166.5
$ws.Range("I67").Value = 3999.8
$ws.Range("K67").Value = 3999.8
$ws.Range("M67").Value = -3141.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1166
$ws.Range("I2").Value = 1166
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1166
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1053
$ws.Range("N2").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H88").Value = 2499.5
$ws.Range("I88").Value = 1999
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 1999
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -1593
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 2499.5
$ws.Range("I91").Value = 1999
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 1999
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -595
$ws.Range("N91").Value = -5808
$ws.Range("H116").Value = 1166
$ws.Range("I116").Value = 1166
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1166
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1128
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1166
$ws.Range("I3").Value = 1166
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1166
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1052
$ws.Range("N3").ClearContents()
$ws.Range("H22").Value = 2779.4
$ws.Range("I22").Value = 2779.4
$ws.Range("K22").Value = 2779.4
$ws.Range("M22").Value = -2606.4
$ws.Range("H29").Value = 949.5
$ws.Range("I29").Value = 949.5
$ws.Range("K29").Value = 949.5
$ws.Range("M29").Value = -660.5
$ws.Range("H86").Value = 5999
$ws.Range("I86").Value = 5999
$ws.Range("K86").Value = 5999
$ws.Range("M86").Value = -4876
$ws.Range("H89").Value = 5999
$ws.Range("I89").Value = 5999
$ws.Range("K89").Value = 29995
$ws.Range("M89").Value = -24379

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6261.615
$ws.Range("I58").Value = 2550.25
$ws.Range("J58").Value = 12199.8
$ws.Range("K58").Value = 2550.25
$ws.Range("L58").Value = 12199.8
$ws.Range("M58").Value = -2347.25
$ws.Range("N58").Value = -12605.8
$ws.Range("H99").Value = 9910
$ws.Range("I99").Value = 9992
$ws.Range("K99").Value = 9992
$ws.Range("M99").Value = -8494
$ws.Range("H116").Value = 39742
$ws.Range("J116").Value = 39742
$ws.Range("L116").Value = 39742
$ws.Range("N116").Value = -48920
$ws.Range("H126").Value = 9910
$ws.Range("I126").Value = 9992
$ws.Range("K126").Value = 29976
$ws.Range("M126").Value = -27506
$ws.Range("H136").Value = 6261.615
$ws.Range("I136").Value = 2550.25
$ws.Range("J136").Value = 12199.8
$ws.Range("K136").Value = 7650.75
$ws.Range("L136").Value = 36599.39999999999
$ws.Range("M136").Value = -5100.75
$ws.Range("N136").Value = -41699.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 66667430
$ws.Range("I4").Value = 820.5714
$ws.Range("K4").Value = 2461.7142
$ws.Range("M4").Value = -2349.7142
$ws.Range("H7").Value = 29.666666
$ws.Range("J7").Value = 35
$ws.Range("L7").Value = 105
$ws.Range("N7").Value = -329
$ws.Range("H102").Value = 1999.3334
$ws.Range("I102").Value = 1999.3334
$ws.Range("K102").Value = 5998.0002
$ws.Range("M102").Value = -3564.0002
$ws.Range("H132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 36.23077
$ws.Range("I2").Value = 26.9
$ws.Range("K2").Value = 26.9
$ws.Range("M2").Value = 86.09999999999999
$ws.Range("H52").Value = 38666.332
$ws.Range("J52").Value = 38666.332
$ws.Range("L52").Value = 38666.332
$ws.Range("N52").Value = -39184.332
$ws.Range("H70").Value = 4985.8
$ws.Range("I70").Value = 4985.8
$ws.Range("K70").Value = 4985.8
$ws.Range("M70").Value = -4715.8
$ws.Range("H73").Value = 4985.8
$ws.Range("I73").Value = 4985.8
$ws.Range("K73").Value = 4985.8
$ws.Range("M73").Value = -4049.8
$ws.Range("H80").Value = 2694.647
$ws.Range("I80").Value = 2673.1333
$ws.Range("K80").Value = 2673.1333
$ws.Range("M80").Value = -1675.1333
$ws.Range("H83").Value = 2694.647
$ws.Range("I83").Value = 2673.1333
$ws.Range("K83").Value = 13365.6665
$ws.Range("M83").Value = -8373.666499999999
$ws.Range("H107").Value = 1074.75
$ws.Range("J107").Value = 1100
$ws.Range("L107").Value = 1100
$ws.Range("N107").Value = -4940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 27000
$ws.Range("J43").Value = 27000
$ws.Range("L43").Value = 27000
$ws.Range("N43").Value = -27386
$ws.Range("H46").Value = 5449
$ws.Range("I46").Value = 5449
$ws.Range("K46").Value = 5449
$ws.Range("M46").Value = -5261
$ws.Range("H82").Value = 896.3333
$ws.Range("I82").Value = 896.3333
$ws.Range("K82").Value = 896.3333
$ws.Range("M82").Value = -535.3333
$ws.Range("H85").Value = 896.3333
$ws.Range("I85").Value = 896.3333
$ws.Range("K85").Value = 896.3333
$ws.Range("M85").Value = 351.6667
$ws.Range("H121").Value = 4055
$ws.Range("I121").Value = 4055
$ws.Range("K121").Value = 4055
$ws.Range("M121").Value = -2308
$ws.Range("H122").Value = 3906
$ws.Range("I122").Value = 3874
$ws.Range("J122").Value = 3916.6667
$ws.Range("K122").Value = 11622
$ws.Range("L122").Value = 11750.0001
$ws.Range("M122").Value = -9172
$ws.Range("N122").Value = -16650.0001
$ws.Range("H136").Value = 10400
$ws.Range("I136").Value = 6000
$ws.Range("K136").Value = 18000
$ws.Range("M136").Value = -15450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 32013
$ws.Range("I22").Value = 32013
$ws.Range("K22").Value = 32013
$ws.Range("M22").Value = -31720
$ws.Range("H81").Value = 1700
$ws.Range("I81").Value = 1700
$ws.Range("K81").Value = 3400
$ws.Range("M81").Value = -2339
$ws.Range("H84").Value = 1700
$ws.Range("I84").Value = 1700
$ws.Range("K84").Value = 17000
$ws.Range("M84").Value = -11696
$ws.Range("H126").Value = 5330
$ws.Range("I126").Value = 5330
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15990
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -13520
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4207.3335
$ws.Range("I132").Value = 2203.9285
$ws.Range("K132").Value = 6611.7855
$ws.Range("M132").Value = -4081.7855

